$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the typo in the ID for the assembly row: EN_A0005 -> EN_A0500
$ws.Range("G2").Value = "EN_A0500"

# Update the active selection to reflect where the user last clicked
$ws.Range("G3").Select()
